$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column E
$ws.Range("E1").Value = "along"

# New constant value for column E data rows
$cutAid = "cut_aid_in_programTRUE"

# Updated numeric values for columns B, C, D (mean, CI_low, CI_high) per row,
# and the new column E value for each data row (rows 2-12).
$data = @(
    @{ B = -0.0394878232656016;  C = -0.0535787331655111;  D = -0.0253969133656922 },
    @{ B = -0.0354902715916673;  C = -0.0561453938825939;  D = -0.0148351493007406 },
    @{ B = -0.00133888279740772; C = -0.0547185798377679;  D = 0.0520408142429524 },
    @{ B = -0.034362429933459;   C = -0.0801901455903655;  D = 0.0114652857234475 },
    @{ B = -0.121455849923609;   C = -0.171665703401423;   D = -0.0712459964457963 },
    @{ B = -0.0485188794061111;  C = -0.112324706905562;   D = 0.0152869480933394 },
    @{ B = -0.0997202653908376;  C = -0.159279311720092;   D = -0.0401612190615832 },
    @{ B = -0.0170813422574052;  C = -0.067562776421911;   D = 0.0334000919071005 },
    @{ B = 0.112828423953117;    C = 0.0431073154165781;   D = 0.182549532489656 },
    @{ B = -0.00814719386724725; C = -0.038439369729519;   D = 0.0221449819950245 },
    @{ B = -0.0675644768188323;  C = -0.0924854230821806;  D = -0.0426435305554841 }
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i].B
    $ws.Cells.Item($row, 3).Value = $data[$i].C
    $ws.Cells.Item($row, 4).Value = $data[$i].D
    $ws.Cells.Item($row, 5).Value = $cutAid
}
